# Auto-generated edit script applying the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.207.90"
$ws.Range("E2").Value = "  +0.95%  "
$ws.Range("D3").Value = "2.267.67"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'305.06"
$ws.Range("E5").Value = "  +0.50%  "
$ws.Range("D6").Value = "'96.21"
$ws.Range("E6").Value = "  +3.86%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.490"
$ws.Range("E9").Value = "  +1.18%  "
$ws.Range("D10").Value = "'35.40"
$ws.Range("E10").Value = "  +9.24%  "
$ws.Range("D11").Value = "'0.0794"
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("D13").Value = "'6.64"
$ws.Range("D14").Value = "2.621.23"
$ws.Range("E14").Value = "  +0.21%  "
$ws.Range("D15").Value = "'14.41"
$ws.Range("E15").Value = "  +0.86%  "
$ws.Range("D16").Value = "2.276.03"
$ws.Range("E16").Value = "  +0.64%  "
$ws.Range("E17").Value = "  +1.31%  "
$ws.Range("D18").Value = "42.131.36"
$ws.Range("E18").Value = "  +0.94%  "
$ws.Range("E19").Value = "  -1.86%  "
$ws.Range("D20").Value = "0.0₃0907"
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("D21").Value = "'5.99"
$ws.Range("D22").Value = "'67.82"
$ws.Range("D23").Value = "'238.22"
$ws.Range("E23").Value = "  -2.56%  "
$ws.Range("D24").Value = "'2.58"
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("E25").Value = "  +0.39%  "
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("D27").Value = "'23.71"
$ws.Range("E27").Value = "  -1.36%  "
$ws.Range("D28").Value = "'37.10"
$ws.Range("E28").Value = "  +5.91%  "
$ws.Range("D29").Value = "'9.53"
$ws.Range("E29").Value = "  -0.54%  "
$ws.Range("E30").Value = "  +1.87%  "
$ws.Range("E31").Value = "  -0.26%  "
$ws.Range("E32").Value = "  -0.36%  "
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  +0.10%  "
$ws.Range("D34").Value = "'3.18"
$ws.Range("E34").Value = "  +5.46%  "
$ws.Range("E35").Value = "  -0.47%  "
$ws.Range("D36").Value = "'17.21"
$ws.Range("E36").Value = "  +1.92%  "
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("E38").Value = "  -0.53%  "
$ws.Range("D39").Value = "'1.83"
$ws.Range("E39").Value = "  +1.98%  "
$ws.Range("E40").Value = "  -1.53%  "
$ws.Range("D41").Value = "'4.07"
$ws.Range("E41").Value = "  +4.01%  "
$ws.Range("E42").Value = "  +9.52%  "
$ws.Range("D43").Value = "1.991.09"
$ws.Range("E43").Value = "  -0.98%  "
$ws.Range("D44").Value = "'18.96"
$ws.Range("E44").Value = "  -4.32%  "
$ws.Range("E45").Value = "  +0.55%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "'9.93"
$ws.Range("E46").Value = "  -3.95%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "'2.92"
$ws.Range("E47").Value = "  +1.13%  "
$ws.Range("D48").Value = "'53.14"
$ws.Range("E48").Value = "  +0.43%  "
$ws.Range("E49").Value = "  +0.40%  "
$ws.Range("E50").Value = "  -1.30%  "
$ws.Range("D51").Value = "'91.23"
$ws.Range("E51").Value = "  -0.46%  "

# Strip the implicit quote-prefix text style picked up above so cells stay unstyled
$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D9").ClearFormats()
$ws.Range("D10").ClearFormats()
$ws.Range("D11").ClearFormats()
$ws.Range("D13").ClearFormats()
$ws.Range("D15").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("D23").ClearFormats()
$ws.Range("D24").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D28").ClearFormats()
$ws.Range("D29").ClearFormats()
$ws.Range("D33").ClearFormats()
$ws.Range("D34").ClearFormats()
$ws.Range("D36").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D41").ClearFormats()
$ws.Range("D44").ClearFormats()
$ws.Range("D46").ClearFormats()
$ws.Range("D47").ClearFormats()
$ws.Range("D48").ClearFormats()
$ws.Range("D51").ClearFormats()
